$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data
$ws.Range("A7").Value = "Cara Core Informática"
$ws.Range("B7").Value = "Programador Python"

# Turn the new email address into a mailto hyperlink (Excel auto-applies
# the built-in "Hyperlink" style to the cell)
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:suporte@caracore.com.br", "", "", "suporte@caracore.com.br")

# Resize the columns to fit their (now longer) contents
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Leave the selection on the cell that was just edited
[void]$ws.Range("C7").Select()
